$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled name used in shared strings ("Prretika Shetty" -> "Preetika Shetty")
# and insert her as a new row 3, shifting the rest of the roster down by one,
# while also updating the ID numbers in column A and scores in column C.
# Remove the now-unused "a" string / G6 helper cell entirely.

$ws.Range("G6").ClearContents()

$ws.Range("A1").Value2 = 52501
$ws.Range("B1").Value = "Jagannath Pidaparthy"
$ws.Range("C1").Value2 = 70

$ws.Range("A2").Value2 = 52502
$ws.Range("B2").Value = "Vishal Patil"
$ws.Range("C2").Value2 = 82

$ws.Range("A3").Value2 = 52503
$ws.Range("B3").Value = "Preetika Shetty"
$ws.Range("C3").Value2 = 98

$ws.Range("A4").Value2 = 52504
$ws.Range("B4").Value = "Sagar Mishra"
$ws.Range("C4").Value2 = 96

$ws.Range("A5").Value2 = 52505
$ws.Range("B5").Value = "Shubham Mishra"
$ws.Range("C5").Value2 = 93

$ws.Range("A6").Value2 = 52506
$ws.Range("B6").Value = "Kanchan Soni"
$ws.Range("C6").Value2 = 95

$ws.Range("A7").Value2 = 52507
$ws.Range("B7").Value = "Jai Lohani"
$ws.Range("C7").Value2 = 92

$ws.Range("A8").Value2 = 52508
$ws.Range("B8").Value = "Korol Dhanda"
$ws.Range("C8").Value2 = 85

$ws.Range("A9").Value2 = 52509
$ws.Range("B9").Value = "kaustubh Srivastava"
$ws.Range("C9").Value2 = 97

$ws.Range("A10").Value2 = 52510
$ws.Range("B10").Value = "Purva Shinde"
$ws.Range("C10").Value2 = 80

$ws.Range("B3").Select()
